$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Duplicate the WeaponData sheet (after itself) so the new sheet inherits
# identical formatting (column widths, number formats, phonetic settings, etc.)
$ws1.Copy($null, $ws1)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "ElementWeaponData"

# Overwrite the copied data with the element-weapon values
$data = @(
    @("SWORD", 3, 0, 5, 0),
    @("LANCE", 2, 0, 0, 5.5),
    @("BOW", 3, 0, 0, 0),
    @("BRASTER", 3.5, 6, 0, 0)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    for ($c = 0; $c -lt 5; $c++) {
        $ws2.Cells.Item($r + 2, $c + 1).Value = $data[$r][$c]
    }
}

# Selections: WeaponData gets the whole-column selection it ends up with after
# the new sheet steals focus; ElementWeaponData keeps the cell last touched.
$null = $ws1.Range("A1:E1048576").Select()
$null = $ws2.Range("B5").Select()

$null = $ws2.Activate()
